$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.1509980484562094
$ws.Range("C2").Value = 1.195282167228728
$ws.Range("D2").Value = 7.57942700223539
$ws.Range("E2").Value = 2.753075916540514
$ws.Range("F2").Value = 2.775751383040755
$ws.Range("G2").Value = 52
$ws.Range("B3").Value = 0.03079452806175957
$ws.Range("C3").Value = 1.189334643345628
$ws.Range("D3").Value = 5.574265677164387
$ws.Range("E3").Value = 2.360988283995578
$ws.Range("F3").Value = 2.384278448858261
$ws.Range("G3").Value = 51
$ws.Range("B4").Value = -0.101817367501059
$ws.Range("C4").Value = 1.016702669087111
$ws.Range("D4").Value = 4.454714751352228
$ws.Range("E4").Value = 2.11061951837659
$ws.Range("F4").Value = 2.129565437855708
$ws.Range("G4").Value = 50
$ws.Range("B5").Value = 0.05611136314148432
$ws.Range("C5").Value = 1.148195161250686
$ws.Range("D5").Value = 5.204774752689769
$ws.Range("E5").Value = 2.281397543763421
$ws.Range("F5").Value = 2.304342309971512
$ws.Range("G5").Value = 49
$ws.Range("B6").Value = -0.06628332611628297
$ws.Range("C6").Value = 1.009343296419617
$ws.Range("D6").Value = 4.545192847316543
$ws.Range("E6").Value = 2.131945789019163
$ws.Range("F6").Value = 2.153465142626193
$ws.Range("G6").Value = 48
$ws.Range("B7").Value = -0.01273090096952459
$ws.Range("C7").Value = 1.122250521844307
$ws.Range("D7").Value = 5.511138675828411
$ws.Range("E7").Value = 2.347581452437468
$ws.Range("F7").Value = 2.380847133747632
$ws.Range("G7").Value = 36
$ws.Range("B8").Value = -0.01901641356679982
$ws.Range("C8").Value = 1.14391309009801
$ws.Range("D8").Value = 5.569593947515463
$ws.Range("E8").Value = 2.359998717693605
$ws.Range("F8").Value = 2.394375341132336
$ws.Range("G8").Value = 35
$ws.Range("B9").Value = 0.05598467855657813
$ws.Range("C9").Value = 1.513431136887689
$ws.Range("D9").Value = 9.181581124177706
$ws.Range("E9").Value = 3.030112394644414
$ws.Range("F9").Value = 3.11742770206278
$ws.Range("G9").Value = 18
$ws.Range("B10").Value = -0.7620237227489945
$ws.Range("C10").Value = 1.215958780976282
$ws.Range("D10").Value = 6.693346821907633
$ws.Range("E10").Value = 2.587150328432353
$ws.Range("F10").Value = 2.593054826775349
$ws.Range("G10").Value = 11
$ws.Range("B11").Value = -0.174551384566527
$ws.Range("C11").Value = 0.5019007433328729
$ws.Range("D11").Value = 0.3135597735887778
$ws.Range("E11").Value = 0.5599640824095575
$ws.Range("F11").Value = 0.5948650978737597
$ws.Range("G11").Value = 5
